$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (mirrors row 2 formatting/content, per the diff) ---

# Style group 2: border box without bottom edge (xf index 2) - no quote prefix
$g2 = @("A3","C3","D3","F3","G3","H3","I3","J3","R3","S3","T3")
foreach ($addr in $g2) {
    $ws.Range($addr).Borders.LineStyle = 1
    $ws.Range($addr).Borders(9).LineStyle = -4142
}
$ws.Range("A3").Value = "DEV"
$ws.Range("C3").Value = "EMEAAD\igrabe"
$ws.Range("D3").Value = "Provide2018"
$ws.Range("F3").Value = "Description PR1"
$ws.Range("G3").Value = "Product 1"
$ws.Range("H3").Value = "Guarding and security Event"
$ws.Range("I3").Value = "S00000912001"
$ws.Range("J3").Value = "30JDF"
$ws.Range("R3").Value = "Bruxelles Triomphe"
$ws.Range("S3").Value = "Sopra Steria Benelux SA/NV"
$ws.Range("T3").Value = "CS-FVAR"

# Style group 1: full border box (xf index 1) - no quote prefix
$g1 = @("N3","X3","Y3")
foreach ($addr in $g1) {
    $ws.Range($addr).Borders.LineStyle = 1
}
$ws.Range("N3").Value = "Purchase requisition"
$ws.Range("X3").Value = "EMEAAD\ncourtin"
$ws.Range("Y3").Value = "EMEAAD\aanciaux"

# Style group 3: border box without bottom edge + numFmt 0.00 (xf index 3) - no quote prefix
$g3 = @("M3","O3","P3")
foreach ($addr in $g3) {
    $ws.Range($addr).Borders.LineStyle = 1
    $ws.Range($addr).Borders(9).LineStyle = -4142
    $ws.Range($addr).NumberFormat = "0.00"
}
$ws.Range("M3").Value = "EUR"
$ws.Range("O3").Value = "STANDARD"
$ws.Range("P3").Value = "Standard"

# Style group 4: border box without bottom edge + quote prefix (xf index 4)
$g4 = @("E3","K3","U3","V3","W3")
foreach ($addr in $g4) {
    $ws.Range($addr).Borders.LineStyle = 1
    $ws.Range($addr).Borders(9).LineStyle = -4142
}
$ws.Range("E3").Value = "'362"
$ws.Range("K3").Value = "'2"
$ws.Range("U3").Value = "'A"
$ws.Range("V3").Value = "'0643"
$ws.Range("W3").Value = "'99"

# Style group 7: border box without bottom edge + numFmt 0.00 + quote prefix (xf index 7)
$g7 = @("L3","Q3")
foreach ($addr in $g7) {
    $ws.Range($addr).Borders.LineStyle = 1
    $ws.Range($addr).Borders(9).LineStyle = -4142
    $ws.Range($addr).NumberFormat = "0.00"
}
$ws.Range("L3").Value = "'2564"
$ws.Range("Q3").Value = "'"

# Style group 8: full border box + quote prefix (xf index 8)
$ws.Range("Z3").Borders.LineStyle = 1
$ws.Range("Z3").Value = "'"

# Style group 6 (Hyperlink cell style): add hyperlink first, then border last
$hl = $ws.Hyperlinks.Add($ws.Range("B3"), "https://voflusoprasttest.p2p.basware.com/edge")
$ws.Range("B3").Borders.LineStyle = 1

Write-Host "Row 3 populated"
